$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Mikal Bridges", "SG,SF,PF", "New York Knicks"),
    @("Buddy Hield", "SG,SF", "Golden State Warriors"),
    @("Josh Giddey", "PG,SG,SF", "Chicago Bulls"),
    @("Julian Champagnie", "SF,PF", "San Antonio Spurs"),
    @("Michael Porter Jr.", "SF,PF", "Denver Nuggets"),
    @("Nikola Vucevic", "PF,C", "Chicago Bulls"),
    @("Bobby Portis", "PF,C", "Milwaukee Bucks"),
    @("Donovan Clingan", "C", "Portland Trail Blazers"),
    @("Brook Lopez", "C", "Milwaukee Bucks"),
    @("DeMar DeRozan", "SF,PF", "Sacramento Kings"),
    @("Tyler Herro", "PG,SG", "Miami Heat"),
    @("De'Aaron Fox", "PG", "Sacramento Kings"),
    @("Bennedict Mathurin", "SG,SF", "Indiana Pacers"),
    @("Evan Mobley", "PF,C", "Cleveland Cavaliers"),
    @("Scottie Barnes", "SG,SF,PF", "Toronto Raptors")
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
